# Applies the recalculated binary "no rank decision" ranking data:
# for each row, updates index/prolificid/name/gender/matrices(range)/race
# to reflect the recomputed within-gender ranking (mat_rank, column J, and
# mat_range bucket, column K, are unaffected by this recalculation).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ Row=2; D=2; E="5c5882fc5bfe7600011197cb"; F="Colleen"; G="female"; H=13.42119510329043; I="White" },
    @{ Row=3; D=3; E="60bd88b8fc436774352f53b9"; F="Annes"; G="female"; H=13.17773416771519; I="Asian" },
    @{ Row=4; D=22; E="608b14a312c099ac00b721b6"; F="Khushi"; G="female"; H=8.277947983434146; I="Asian" },
    @{ Row=5; D=21; E="5c0e89c6c323400001e6c4a5"; F="Bri"; G="female"; H=8.218874334828817; I="Black or African American" },
    @{ Row=6; D=19; E="60b45e9961dd412bfb6780f8"; F="Jewel"; G="female"; H=8.21192345112825; I="Black or African American" },
    @{ Row=7; D=33; E="60cb36ee9f58331a33cf5506"; F="Shaniek"; G="female"; H=5.441970684512863; I="Black or African American" },
    @{ Row=8; D=32; E="6036f9b3b1842f8b659b18c7"; F="Kellie"; G="female"; H=5.381459162249058; I="White" },
    @{ Row=9; D=30; E="60d5775a99b502eec8cf56b4"; F="Shadaisia"; G="female"; H=5.321845954194636; I="Black or African American" },
    @{ Row=10; D=34; E="5e96194b0a9fe909389e9f7b"; F="Tina"; G="female"; H=4.498467056693604; I="White" },
    @{ Row=11; D=35; E="6077db0613ce87b4a62a78f9"; F="Lori"; G="female"; H=4.222996349665409; I="White" },
    @{ Row=12; D=41; E="60bfcf5805c5ae12a546f9f3"; F="Giana"; G="female"; H=2.390791975163696; I="White" },
    @{ Row=13; D=44; E="60c0e5899d387663c07eb3a4"; F="Nansi"; G="female"; H=1.089220531548616; I="Asian" },
    @{ Row=14; D=2; E="5e2522d6b734b47915f88275"; F="Corey"; G="male"; H=14.11239547175637; I="White" },
    @{ Row=15; D=3; E="601d69a993d94008fb2b25dc"; F="Quinterius"; G="male"; H=13.09487473480318; I="Black or African American" },
    @{ Row=16; D=22; E="60db4fde6193c50664c9c478"; F="Edosagbe"; G="male"; H=8.22111200880744; I="Black or African American" },
    @{ Row=17; D=26; E="5dd671942b033b5ec8bc97b4"; F="Juan"; G="male"; H=7.429121582096163; I="Hispanic" },
    @{ Row=18; D=27; E="5ff8ad350d084e10f500e48a"; F="Drew"; G="male"; H=6.324528075904071; I="White" },
    @{ Row=19; D=29; E="60b83826821417f8e484a207"; F="Eli"; G="male"; H=6.243826188088984; I="White" },
    @{ Row=20; D=32; E="60bf9943e4e04642d4634ecc"; F="Jamarii"; G="male"; H=5.27722767756892; I="Black or African American" },
    @{ Row=21; D=33; E="60b322994d0b901954690036"; F="Brennan"; G="male"; H=5.186042016282854; I="White" },
    @{ Row=22; D=30; E="60c2341fe95d71ee52c043f0"; F="Matthew"; G="male"; H=5.141087836715284; I="White" },
    @{ Row=23; D=44; E="60b091ed11ccda59e3fc7761"; F="Myles"; G="male"; H=3.417079858592328; I="Black or African American" },
    @{ Row=24; D=49; E="6088fc724afd5c008db33e9d"; F="Masuf"; G="male"; H=1.153463192899035; I="Asian" },
    @{ Row=25; D=50; E="6097b95056caf5ebb2720002"; F="Damian"; G="male"; H=0.3269558257719956; I="Black or African American" }
)

foreach ($r in $rowData) {
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("F" + $r.Row).Value = $r.F
    $ws.Range("G" + $r.Row).Value = $r.G
    $ws.Range("H" + $r.Row).Value = $r.H
    $ws.Range("I" + $r.Row).Value = $r.I
}
